$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: D2 changes from "sCs" to "FAPs", and several numeric values change
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2197353333333333
$ws.Range("H2").Value = 0.659206
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07823633333333334
$ws.Range("N2").Value = 0.234709
$ws.Range("O2").Value = 0.07325462794193288
$ws.Range("P2").Value = 0.07325462794193287
$ws.Range("Q2").Value = 0.01719128678377778
$ws.Range("R2").Value = 0.154721581054
$ws.Range("S2").Value = 0.07325462794193288
$ws.Range("T2").Value = 0.07325462794193287

# Add new row 3
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Nlgn1"
$ws.Range("C3").Value = "Nrxn1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2197353333333333
$ws.Range("H3").Value = 0.659206
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.989769
$ws.Range("N3").Value = 2.969307
$ws.Range("O3").Value = 0.9267453720580672
$ws.Range("P3").Value = 0.9267453720580671
$ws.Range("Q3").Value = 0.217487221138
$ws.Range("R3").Value = 1.957384990242
$ws.Range("S3").Value = 0.9267453720580672
$ws.Range("T3").Value = 0.9267453720580671
